$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ between rows 4 and 5 (per the target diff);
# columns with identical values in both rows (C,S,T,U,V,W,Y,AA,AD,AE,AG,AW,AX, etc.)
# are intentionally left untouched so their XML representation is not disturbed.
$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB")

foreach ($col in $cols) {
    $addr4 = "$col" + "4"
    $addr5 = "$col" + "5"
    $val4 = $ws.Range($addr4).Value2
    $val5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value = $val5
    $ws.Range($addr5).Value = $val4
}
